$wb = $excel.ActiveWorkbook

# --- Update the view/selection on the existing sheet ("page-1_table-1") ---
$ws1 = $wb.Worksheets.Item(1)
[void]$ws1.Activate()
[void]$ws1.Range("S1:S14").Select()
$excel.ActiveWindow.ScrollColumn = 17

# --- Add the new "Irrigation" worksheet after the existing last sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Irrigation"

# Header row (bold, matches style index 1 used by the source data)
$ws.Range("B4").Value = "Region"
$ws.Range("C4").Value = 2001
$ws.Range("D4").Value = 2023
$ws.Range("B4:D4").Font.Bold = $true

# Region / 2001 data; column D is 2001 value * 2
$regions = @(
    "Boucle du Mouhoun",
    "Cascades",
    "Centre",
    "Centre-Est",
    "Centre-Nord",
    "Centre-Ouest",
    "Centre-Sud",
    "Est",
    "Hauts-Bassins",
    "Nord",
    "Plateau Central",
    "Sahel",
    "Sud-Ouest"
)
$values2001 = @(4391, $null, 1495, 2210, 1135, 555, 2210, 580, 13165, 875, 1700, 270, 325)

for ($i = 0; $i -lt $regions.Length; $i++) {
    $row = 5 + $i
    $ws.Range("B$row").Value = $regions[$i]
    if ($null -ne $values2001[$i]) {
        $ws.Range("C$row").Value = $values2001[$i]
    }
    $ws.Range("D$row").Formula = "=C$row*2"
}

$ws.Columns.Item(2).ColumnWidth = 16.9167

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

[void]$ws.Range("F8").Select()
